$wb = $excel.ActiveWorkbook
$wsJs = $wb.Worksheets.Item("flair.js")
$wsJs.Range("D151:D200").Value = '": "'
$wsJs.Activate()
$wsJs.Range("H133:H198").Select()
